$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.25
$ws.Range("I2").Value = 3.4
$ws.Range("J2").Value = 3
$ws.Range("L2").Value = 4
$ws.Range("X2").Value = 10
$ws.Range("AV2").Value = 67
$ws.Range("AZ2").Value = 67
$ws.Range("BA2").Value = 101
